$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it
# stored as TEXT (matches the source data, which is all inline/shared
# strings, not numbers). Temporarily force a Text number format, assign
# the value, then restore the cell to the unstyled "Normal" style so the
# saved XML carries no stray style index (matching the original file).
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row-by-row value updates (price column D, volume(1h) column E)
Set-TextValue "D2" "69.345.42"
$ws.Range("E2").Value = "  -1.31%  "

Set-TextValue "D3" "3.536.82"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue "D5" "573.38"
$ws.Range("E5").Value = "  -0.64%  "

Set-TextValue "D6" "182.90"
$ws.Range("E6").Value = "  -4.21%  "

Set-TextValue "D7" "3.527.67"
$ws.Range("E7").Value = "  -1.68%  "

Set-TextValue "D8" "0.616"
$ws.Range("E8").Value = "  -2.74%  "

$ws.Range("E9").Value = "  +0.15%  "

Set-TextValue "D10" "0.189"
$ws.Range("E10").Value = "  +6.29%  "

Set-TextValue "D11" "0.641"
$ws.Range("E11").Value = "  -3.05%  "

Set-TextValue "D12" "53.83"
$ws.Range("E12").Value = "  -5.00%  "

Set-TextValue "D13" "0.0000302"
$ws.Range("E13").Value = "  +0.56%  "

Set-TextValue "D14" "9.50"
$ws.Range("E14").Value = "  -3.06%  "

Set-TextValue "D15" "4.103.30"
$ws.Range("E15").Value = "  -1.22%  "

Set-TextValue "D16" "19.35"
$ws.Range("E16").Value = "  -4.02%  "

Set-TextValue "D17" "3.530.19"
$ws.Range("E17").Value = "  -1.37%  "

Set-TextValue "D18" "69.275.36"
$ws.Range("E18").Value = "  -1.03%  "

Set-TextValue "D19" "12.57"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("E20").Value = "  -1.19%  "

Set-TextValue "D21" "534.74"
$ws.Range("E21").Value = "  +12.81%  "

Set-TextValue "D22" "1.04"
$ws.Range("E22").Value = "  -0.68%  "

Set-TextValue "D23" "19.74"
$ws.Range("E23").Value = "  +0.22%  "

Set-TextValue "D24" "4.97"
$ws.Range("E24").Value = "  -3.09%  "

Set-TextValue "D25" "4.38"
$ws.Range("E25").Value = "  +0.55%  "

Set-TextValue "D26" "94.37"
$ws.Range("E26").Value = "  +6.44%  "

Set-TextValue "D27" "11.12"
$ws.Range("E27").Value = "  -0.16%  "

Set-TextValue "D28" "2.93"
$ws.Range("E28").Value = "  -4.37%  "

Set-TextValue "D29" "9.10"
$ws.Range("E29").Value = "  -2.02%  "

Set-TextValue "D30" "31.87"

Set-TextValue "D31" "7.37"
$ws.Range("E31").Value = "  -5.74%  "

Set-TextValue "D32" "12.59"
$ws.Range("E32").Value = "  +3.90%  "

Set-TextValue "D33" "64.95"
$ws.Range("E33").Value = "  -1.95%  "

Set-TextValue "D34" "0.115"
$ws.Range("E34").Value = "  -4.36%  "

Set-TextValue "D35" "572.74"
$ws.Range("E35").Value = "  -2.77%  "

Set-TextValue "D36" "3.16"
$ws.Range("E36").Value = "  +7.36%  "

Set-TextValue "D37" "38.49"
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("E38").Value = "  +0.22%  "

Set-TextValue "D39" "0.401"
$ws.Range("E39").Value = "  +0.30%  "

Set-TextValue "D40" "0.0₃0766"
$ws.Range("E40").Value = "  -5.18%  "

Set-TextValue "D43" "3.37"
$ws.Range("E43").Value = "  -4.70%  "

$ws.Range("E44").Value = "  +4.13%  "

Set-TextValue "D45" "2.97"
$ws.Range("E45").Value = "  -4.38%  "

Set-TextValue "D46" "0.0444"
$ws.Range("E46").Value = "  -0.54%  "

Set-TextValue "D47" "3.167.61"
$ws.Range("E47").Value = "  -2.10%  "

Set-TextValue "D48" "9.23"
$ws.Range("E48").Value = "  -4.12%  "

$ws.Range("E49").Value = "  -1.91%  "

Set-TextValue "D50" "0.997"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("E51").Value = "  -2.56%  "

# Rows 41 and 42 swapped their Coin/Link/Price content, with new Volume(1h) values
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D41" "3.09"
$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.133"
$ws.Range("E42").Value = "  -6.85%  "

